# Update horarios workbook (Linea 141) with the latest scrape data.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Sheet: LP1912 ---
$ws1.Cells.Item(2, 1).Value = "Última actualización: 07:28:14"
$ws1.Cells.Item(3, 1).Value = "Total filas: 91"
$ws1.Cells.Item(39, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(40, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(47, 3).Value = "14_ABASTO"
$ws1.Cells.Item(48, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(52, 1).Value = "05:52:07"
$ws1.Cells.Item(52, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(52, 4).Value = 73
$ws1.Cells.Item(53, 1).Value = "06:59:37"
$ws1.Cells.Item(53, 3).Value = "15_ABASTO"
$ws1.Cells.Item(53, 4).Value = 6
$ws1.Cells.Item(58, 1).Value = "06:59:37"
$ws1.Cells.Item(58, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(58, 4).Value = 17
$ws1.Cells.Item(59, 1).Value = "05:52:07"
$ws1.Cells.Item(59, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(59, 4).Value = 84
$ws1.Cells.Item(62, 1).Value = "07:28:14"
$ws1.Cells.Item(62, 2).Value = "07:30"
$ws1.Cells.Item(62, 4).Value = 2
$ws1.Cells.Item(64, 2).Value = "07:31"
$ws1.Cells.Item(64, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(64, 4).Value = 32
$ws1.Cells.Item(65, 1).Value = "07:28:14"
$ws1.Cells.Item(65, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(65, 4).Value = 4
$ws1.Cells.Item(66, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(67, 1).Value = "05:52:07"
$ws1.Cells.Item(67, 2).Value = "07:32"
$ws1.Cells.Item(67, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(67, 4).Value = 100
$ws1.Cells.Item(68, 1).Value = "07:28:14"
$ws1.Cells.Item(68, 2).Value = "07:35"
$ws1.Cells.Item(68, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(68, 4).Value = 7
$ws1.Cells.Item(69, 2).Value = "07:36"
$ws1.Cells.Item(69, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(69, 4).Value = 37
$ws1.Cells.Item(70, 1).Value = "07:28:14"
$ws1.Cells.Item(70, 2).Value = "07:37"
$ws1.Cells.Item(70, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(70, 4).Value = 9
$ws1.Cells.Item(71, 1).Value = "07:28:14"
$ws1.Cells.Item(71, 2).Value = "07:39"
$ws1.Cells.Item(71, 3).Value = "10_OLMOS"
$ws1.Cells.Item(71, 4).Value = 11
$ws1.Cells.Item(72, 2).Value = "07:47"
$ws1.Cells.Item(72, 3).Value = "14_ABASTO"
$ws1.Cells.Item(72, 4).Value = 48
$ws1.Cells.Item(73, 1).Value = "07:28:14"
$ws1.Cells.Item(73, 2).Value = "07:47"
$ws1.Cells.Item(73, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(73, 4).Value = 19
$ws1.Cells.Item(74, 1).Value = "07:28:14"
$ws1.Cells.Item(74, 2).Value = "07:48"
$ws1.Cells.Item(74, 3).Value = "14_ABASTO"
$ws1.Cells.Item(74, 4).Value = 20
$ws1.Cells.Item(75, 1).Value = "07:28:14"
$ws1.Cells.Item(75, 2).Value = "07:51"
$ws1.Cells.Item(75, 3).Value = "215D_EL PATO"
$ws1.Cells.Item(75, 4).Value = 23
$ws1.Cells.Item(76, 1).Value = "07:28:14"
$ws1.Cells.Item(76, 2).Value = "07:55"
$ws1.Cells.Item(76, 3).Value = "10_OLMOS"
$ws1.Cells.Item(76, 4).Value = 27
$ws1.Cells.Item(77, 1).Value = "07:28:14"
$ws1.Cells.Item(77, 2).Value = "08:00"
$ws1.Cells.Item(77, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(77, 4).Value = 32
$ws1.Cells.Item(78, 1).Value = "06:21:22"
$ws1.Cells.Item(78, 2).Value = "08:01"
$ws1.Cells.Item(78, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(78, 4).Value = 100
$ws1.Cells.Item(79, 1).Value = "07:28:14"
$ws1.Cells.Item(79, 2).Value = "08:03"
$ws1.Cells.Item(79, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(79, 4).Value = 35
$ws1.Cells.Item(80, 2).Value = "08:06"
$ws1.Cells.Item(80, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(80, 4).Value = 67
$ws1.Cells.Item(81, 1).Value = "07:28:14"
$ws1.Cells.Item(81, 2).Value = "08:11"
$ws1.Cells.Item(81, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(81, 4).Value = 43
$ws1.Cells.Item(82, 1).Value = "07:28:14"
$ws1.Cells.Item(82, 2).Value = "08:12"
$ws1.Cells.Item(82, 3).Value = "15_ABASTO"
$ws1.Cells.Item(82, 4).Value = 44
$ws1.Cells.Item(82, 5).Value = "LP1912"
$ws1.Cells.Item(83, 1).Value = "07:28:14"
$ws1.Cells.Item(83, 2).Value = "08:21"
$ws1.Cells.Item(83, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(83, 4).Value = 53
$ws1.Cells.Item(83, 5).Value = "LP1912"
$ws1.Cells.Item(84, 1).Value = "06:59:37"
$ws1.Cells.Item(84, 2).Value = "08:22"
$ws1.Cells.Item(84, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(84, 4).Value = 83
$ws1.Cells.Item(84, 5).Value = "LP1912"
$ws1.Cells.Item(85, 1).Value = "07:28:14"
$ws1.Cells.Item(85, 2).Value = "08:23"
$ws1.Cells.Item(85, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(85, 4).Value = 55
$ws1.Cells.Item(85, 5).Value = "LP1912"
$ws1.Cells.Item(86, 1).Value = "07:28:14"
$ws1.Cells.Item(86, 2).Value = "08:23"
$ws1.Cells.Item(86, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(86, 4).Value = 55
$ws1.Cells.Item(86, 5).Value = "LP1912"
$ws1.Cells.Item(87, 1).Value = "07:28:14"
$ws1.Cells.Item(87, 2).Value = "08:27"
$ws1.Cells.Item(87, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(87, 4).Value = 59
$ws1.Cells.Item(87, 5).Value = "LP1912"
$ws1.Cells.Item(88, 1).Value = "07:28:14"
$ws1.Cells.Item(88, 2).Value = "08:42"
$ws1.Cells.Item(88, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(88, 4).Value = 74
$ws1.Cells.Item(88, 5).Value = "LP1912"
$ws1.Cells.Item(89, 1).Value = "07:28:14"
$ws1.Cells.Item(89, 2).Value = "08:44"
$ws1.Cells.Item(89, 3).Value = "14_ABASTO"
$ws1.Cells.Item(89, 4).Value = 76
$ws1.Cells.Item(89, 5).Value = "LP1912"
$ws1.Cells.Item(90, 1).Value = "07:28:14"
$ws1.Cells.Item(90, 2).Value = "08:54"
$ws1.Cells.Item(90, 3).Value = "17_ROMERO"
$ws1.Cells.Item(90, 4).Value = 86
$ws1.Cells.Item(90, 5).Value = "LP1912"
$ws1.Cells.Item(91, 1).Value = "07:28:14"
$ws1.Cells.Item(91, 2).Value = "09:02"
$ws1.Cells.Item(91, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(91, 4).Value = 94
$ws1.Cells.Item(91, 5).Value = "LP1912"
$ws1.Cells.Item(92, 1).Value = "07:28:14"
$ws1.Cells.Item(92, 2).Value = "09:11"
$ws1.Cells.Item(92, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(92, 4).Value = 103
$ws1.Cells.Item(92, 5).Value = "LP1912"
$ws1.Cells.Item(93, 1).Value = "07:28:14"
$ws1.Cells.Item(93, 2).Value = "09:17"
$ws1.Cells.Item(93, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(93, 4).Value = 109
$ws1.Cells.Item(93, 5).Value = "LP1912"
$ws1.Cells.Item(94, 1).Value = "07:28:14"
$ws1.Cells.Item(94, 2).Value = "09:21"
$ws1.Cells.Item(94, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(94, 4).Value = 113
$ws1.Cells.Item(94, 5).Value = "LP1912"
$ws1.Cells.Item(95, 1).Value = "07:28:14"
$ws1.Cells.Item(95, 2).Value = "09:23"
$ws1.Cells.Item(95, 3).Value = "17_ROMERO"
$ws1.Cells.Item(95, 4).Value = 115
$ws1.Cells.Item(95, 5).Value = "LP1912"
$ws1.Cells.Item(96, 1).Value = "07:28:14"
$ws1.Cells.Item(96, 2).Value = "09:24"
$ws1.Cells.Item(96, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(96, 4).Value = 116
$ws1.Cells.Item(96, 5).Value = "LP1912"

# --- Sheet: LP1912-215 ---
$ws2.Cells.Item(2, 1).Value = "Última actualización: 07:28:14"
$ws2.Cells.Item(3, 1).Value = "Total filas: 17"
$ws2.Cells.Item(20, 1).Value = "07:28:14"
$ws2.Cells.Item(20, 4).Value = 23
$ws2.Cells.Item(21, 1).Value = "07:28:14"
$ws2.Cells.Item(21, 4).Value = 55
$ws2.Cells.Item(22, 1).Value = "07:28:14"
$ws2.Cells.Item(22, 2).Value = "09:02"
$ws2.Cells.Item(22, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(22, 4).Value = 94
$ws2.Cells.Item(22, 5).Value = "LP1912"

# --- Sheet: 6203-6173 ---
$ws3.Cells.Item(2, 1).Value = "Última actualización: 07:28:14"
$ws3.Cells.Item(3, 1).Value = "Total filas: 17"
$ws3.Cells.Item(16, 1).Value = "07:28:14"
$ws3.Cells.Item(16, 4).Value = 7
$ws3.Cells.Item(19, 1).Value = "07:28:14"
$ws3.Cells.Item(19, 2).Value = "08:10"
$ws3.Cells.Item(19, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(19, 4).Value = 42
$ws3.Cells.Item(19, 5).Value = "L6203"
$ws3.Cells.Item(20, 1).Value = "06:49:33"
$ws3.Cells.Item(20, 2).Value = "08:33"
$ws3.Cells.Item(20, 4).Value = 104
$ws3.Cells.Item(21, 1).Value = "07:28:14"
$ws3.Cells.Item(21, 2).Value = "08:38"
$ws3.Cells.Item(21, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(21, 4).Value = 70
$ws3.Cells.Item(21, 5).Value = "L6173"
$ws3.Cells.Item(22, 1).Value = "07:28:14"
$ws3.Cells.Item(22, 2).Value = "09:09"
$ws3.Cells.Item(22, 3).Value = "215D_LA PLATA"
$ws3.Cells.Item(22, 4).Value = 101
$ws3.Cells.Item(22, 5).Value = "L6203"
